$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 6 (the old trial data rows), keeping only the header
# row and the single updated data row.
$ws.Range("A3:D6").EntireRow.Delete() | Out-Null

# Update the remaining data row (row 2) with the new trial results.
$ws.Range("A2").Value = "Dispenser"
$ws.Range("B2").Value = "Water"
$ws.Range("C2").Value = 165
$ws.Range("D2").Value = 160
